# Timesheet update: add Jared's Sprint 4 hours + notes, and move the
# active selection to reflect where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data for Sprint 4 / Jared (column E) ---------------------------
# Row 22: hours column (numeric), Row 23: notes column (text)
$ws.Range("E22").Value = 4
$ws.Range("E23").Value = "Finished changes made check, encryption"

# Match the centered alignment used by the rest of the hours/notes columns
# (e.g. D19/C19, which share this same look).
$ws.Range("E22:E23").HorizontalAlignment = -4108

# --- Reflect the final cursor position / scroll state -------------------
$ws.Range("F23").Select()

Write-Output "Sprint 4 / Jared entries added"
